$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Nicholas Marzadro"
$ws.Range("B18").Value = "Samuele Kettamier | SBARX"
$ws.Range("C18").Value = "Leonardo  Parisi  | MediaserT"
$ws.Range("D18").Value = "Mattia Festi | Shark Attack"
$ws.Range("E18").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("F18").Value = "Davide  Bazzano  | iMontagna"
